$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E16").Value = "2305"
$ws.Range("E17").Value = "2304"
$ws.Range("E18").Value = "2303"
$ws.Range("E19").Value = "2302"
